$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Structural edits: insert the new "Sub topics" column and the 3 extra
#    rows needed for the Data-Ingestion sub-bullets.
# ---------------------------------------------------------------------------

# Insert a new column before B (old column B "Resources" shifts to C).
$ws.Columns("B:B").Insert()

# Give the new column B the same width as column A (best effort - the COM
# layer quantises ColumnWidth to whole pixels so this lands as close as the
# host allows to column A's width).
$ws.Columns("B:B").ColumnWidth = $ws.Columns("A:A").ColumnWidth

# Insert 3 new rows right after row 4 ("Data Ingestion & Data Versioning")
# for the sub-topic bullets; everything below shifts down by 3 rows.
$ws.Rows("5:7").Insert()

# ---------------------------------------------------------------------------
# 2. New header cell
# ---------------------------------------------------------------------------
$ws.Range("B1").Value2 = "Sub topics"

# ---------------------------------------------------------------------------
# 3. Sub-topic values in the new column B, rows 4-7
# ---------------------------------------------------------------------------
$ws.Range("B4").Value2 = "TFX for data ingestion from files/services"
$ws.Range("B5").Value2 = "train-test splits"
$ws.Range("B6").Value2 = "combine multiple data exports into one all-encompassing dataset"
$ws.Range("B7").Value2 = "strategies to ingest different forms of data (structured, text, and images)"

# Rows 6-7 get the smaller Times New Roman font used for the longer bullets.
$r6 = $ws.Range("B6")
$r6.Font.Name = "Times New Roman"
$r6.Font.Size = 9
$r6.Font.Color = 3355443
$r6.Font.Family = 1

$r7 = $ws.Range("B7")
$r7.Font.Name = "Times New Roman"
$r7.Font.Size = 9
$r7.Font.Color = 3355443
$r7.Font.Family = 1

# ---------------------------------------------------------------------------
# 4. Boxed border + left/center alignment around A4:A7 (applied to each cell
#    before merging, then the range is merged into one cell).
# ---------------------------------------------------------------------------

# A4 - top of the box: keep left/right/top, drop the bottom edge.
$a4 = $ws.Range("A4")
$a4.Borders.Item(9).LineStyle = -4142
$a4.HorizontalAlignment = -4131
$a4.VerticalAlignment = -4108

# A5 - middle of the box: keep left/right, drop top and bottom.
$a5 = $ws.Range("A5")
$a5.Borders.Item(8).LineStyle = -4142
$a5.Borders.Item(9).LineStyle = -4142
$a5.HorizontalAlignment = -4131
$a5.VerticalAlignment = -4108

# A6 - middle of the box: keep left/right, drop top and bottom.
$a6 = $ws.Range("A6")
$a6.Borders.Item(8).LineStyle = -4142
$a6.Borders.Item(9).LineStyle = -4142
$a6.HorizontalAlignment = -4131
$a6.VerticalAlignment = -4108

# A7 - bottom of the box: keep left/right/bottom, drop the top edge.
$a7 = $ws.Range("A7")
$a7.Borders.Item(8).LineStyle = -4142
$a7.HorizontalAlignment = -4131
$a7.VerticalAlignment = -4108

# Merge the four cells into the single "Data Ingestion & Data Versioning" cell.
$ws.Range("A4:A7").Merge()

# ---------------------------------------------------------------------------
# 5. Selection / active cell, matching the commit's final UI state.
# ---------------------------------------------------------------------------
$ws.Range("A4:A7").Select()

Write-Output "edit complete"
